# Auto-generated edit script applying cell-value changes per the commit diff.
# Each sheet's edits are applied by iterating over (row, col, value) triples.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# ALC: 7 cell updates
$updates = @(
  @(19,8,733.125),
  @(19,9,472),
  @(19,10,994.25),
  @(19,11,472),
  @(19,12,994.25),
  @(19,13,-297),
  @(19,14,-1344.25)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("ARM")
# ARM: 25 cell updates
$updates = @(
  @(32,8,13238.8),
  @(32,9,7014.074),
  @(32,10,26167.076),
  @(32,11,7014.074),
  @(32,12,26167.076),
  @(32,13,-6727.074),
  @(32,14,-26741.076),
  @(74,8,33874.03),
  @(74,9,78636.30499999999),
  @(74,10,3247.2104),
  @(74,11,78636.30499999999),
  @(74,12,3247.2104),
  @(74,13,-77762.30499999999),
  @(74,14,-4995.2104),
  @(77,8,33874.03),
  @(77,9,78636.30499999999),
  @(77,10,3247.2104),
  @(77,11,393181.525),
  @(77,12,16236.052),
  @(77,13,-388813.525),
  @(77,14,-24972.052),
  @(125,8,50721),
  @(125,10,52151.25),
  @(125,12,52151.25),
  @(125,14,-61991.25)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("BSM")
# BSM: 23 cell updates
$updates = @(
  @(99,8,5297648),
  @(99,9,118753.336),
  @(99,11,118753.336),
  @(99,13,-117255.336),
  @(105,8,39028.406),
  @(105,9,48981.047),
  @(105,11,48981.047),
  @(105,13,-47234.047),
  @(107,8,3179.8823),
  @(107,9,2607.3572),
  @(107,10,5851.6665),
  @(107,11,2607.3572),
  @(107,12,5851.6665),
  @(107,13,-687.3571999999999),
  @(107,14,-9691.666499999999),
  @(117,8,79905),
  @(117,10,79905),
  @(117,12,79905),
  @(117,14,-89083),
  @(135,8,48916.75),
  @(135,10,48916.75),
  @(135,12,48916.75),
  @(135,14,-59056.75)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("CRP")
# CRP: 39 cell updates
$updates = @(
  @(31,8,3942.8572),
  @(31,9,2988.375),
  @(31,10,5215.5),
  @(31,11,2988.375),
  @(31,12,5215.5),
  @(31,13,-2693.375),
  @(31,14,-5805.5),
  @(34,8,3942.8572),
  @(34,9,2988.375),
  @(34,10,5215.5),
  @(34,11,2988.375),
  @(34,12,5215.5),
  @(34,13,-2786.375),
  @(34,14,-5619.5),
  @(99,8,3476582.2),
  @(99,9,4623.5),
  @(99,10,10420500),
  @(99,11,4623.5),
  @(99,12,10420500),
  @(99,13,-3125.5),
  @(99,14,-10423496),
  @(105,8,4887.25),
  @(105,9,1033),
  @(105,10,7199.8),
  @(105,11,1033),
  @(105,12,7199.8),
  @(105,13,714),
  @(105,14,-10693.8),
  @(116,8,51387),
  @(116,10,51387),
  @(116,12,51387),
  @(116,14,-60565),
  @(126,8,3476582.2),
  @(126,9,4623.5),
  @(126,10,10420500),
  @(126,11,13870.5),
  @(126,12,31261500),
  @(126,13,-11400.5),
  @(126,14,-31266440)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("CUL")
# CUL: 37 cell updates
$updates = @(
  @(5,8,1497.6428),
  @(5,9,1193.8889),
  @(5,10,2044.4),
  @(5,11,3581.6667),
  @(5,12,6133.200000000001),
  @(5,13,-3469.6667),
  @(5,14,-6357.200000000001),
  @(56,8,5505),
  @(56,9,5505),
  @(56,11,5505),
  @(56,13,-4975),
  @(113,8,91812.27),
  @(113,10,112081.664),
  @(113,12,336244.992),
  @(113,14,-340584.992),
  @(135,8,1497.6428),
  @(135,9,1193.8889),
  @(135,10,2044.4),
  @(135,11,10745.0001),
  @(135,12,18399.6),
  @(135,13,-8210.000099999999),
  @(135,14,-23469.6),
  @(137,8,10679.177),
  @(137,9,6470),
  @(137,10,12433),
  @(137,11,19410),
  @(137,12,37299),
  @(137,13,-14310),
  @(137,14,-47499),
  @(139,8,7581.68),
  @(139,10,8916.947),
  @(139,12,26750.841),
  @(139,14,-37030.841),
  @(140,8,2599.3333),
  @(140,9,1652.909),
  @(140,11,4958.727000000001),
  @(140,13,221.2729999999992)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("GSM")
# GSM: 7 cell updates
$updates = @(
  @(132,8,3914.074),
  @(132,9,3262.2632),
  @(132,10,5462.125),
  @(132,11,9786.7896),
  @(132,12,16386.375),
  @(132,13,-7256.7896),
  @(132,14,-21446.375)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("LTW")
# LTW: 43 cell updates
$updates = @(
  @(7,8,40973.082),
  @(7,10,270000),
  @(7,12,270000),
  @(7,14,-270224),
  @(16,8,1397.85),
  @(16,9,1377.6666),
  @(16,11,1377.6666),
  @(16,13,-1207.6666),
  @(22,8,975.16),
  @(22,9,807.3),
  @(22,11,807.3),
  @(22,13,-512.3),
  @(27,8,975.16),
  @(27,9,807.3),
  @(27,11,807.3),
  @(27,13,-700.3),
  @(46,8,2800),
  @(46,9,600),
  @(46,10,5000),
  @(46,11,600),
  @(46,12,5000),
  @(46,13,-412),
  @(46,14,-5376),
  @(55,8,1285.2264),
  @(55,9,629.8857400000001),
  @(55,11,629.8857400000001),
  @(55,13,-456.8857400000001),
  @(58,8,4635),
  @(58,9,4635),
  @(58,11,4635),
  @(58,13,-4375),
  @(64,8,3500),
  @(64,10,0),
  @(64,12,0),
  @(64,14,$null),
  @(67,8,3500),
  @(67,10,0),
  @(67,12,0),
  @(67,14,$null),
  @(126,8,40973.082),
  @(126,10,270000),
  @(126,12,810000),
  @(126,14,-814940)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets("WVR")
# WVR: 4 cell updates
$updates = @(
  @(20,8,20500),
  @(20,10,0),
  @(20,12,0),
  @(20,14,$null)
)
foreach ($u in $updates) {
  $r = $u[0]; $c = $u[1]; $v = $u[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}
